# H5set_free_list_limits_questions.xlsx - minor copy edits
# - tweak a few question wordings (punctuation / "write and read data")
# - turn word-wrap on for the "notes" column entries that were missing it

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Text / copy edits -------------------------------------------------
$ws.Range("A9").Value  = "If you were going to use this function, is there additional information you would like to know that is not in the entry? If yes, note what that information is."
$ws.Range("A20").Value = "If you were going to use this function, is there additional information you would like to know that is not in the entry? If yes, note what that information is."
$ws.Range("A10").Value = "What type of user would use this function? (HDF library developer,  tool or application developer, user reading an HDF dataset) - answer in notes column."
$ws.Range("A12").Value = "In this section, pretend you are a new user to HDF who is trying to write and read data, and answer these questions with that frame of mind."

# --- Formatting: enable wrap text on the "notes" column (E) ------------
$notesRows = 4,5,6,7,8,9,10,11,14,15,16,17,18,19,20,21
foreach ($r in $notesRows) {
    $ws.Cells.Item($r, 5).WrapText = $true
}

# --- Row heights that Excel re-flowed after the wrap/width tweaks ------
$ws.Rows.Item(5).RowHeight  = 30
$ws.Rows.Item(10).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 30
$ws.Rows.Item(18).RowHeight = 30
$ws.Rows.Item(19).RowHeight = 30

# --- Cosmetic: leave the selection where the author left it on save ----
[void]$ws.Range("G26").Select()
